$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Centraal Station" stop (row 16, A16) to "Centraal Station_A"
$ws.Range("A16").Value = "Centraal Station_A"

# Leave the selection on the edited cell (A16), matching the saved view state
$ws.Range("A16").Select()
